# Generate Report for Handoff
# Refresh the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" timestamps
# for 7c391cd0-1d38-42be-a45a-8f7a7f348bc8.md now that it has been (re)handed off.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G5").Value = "2016-09-01 12:46:35"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H5").Value = "2016-09-01 12:46:30"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H5").Value = "2016-09-01 12:46:35"
